# Edit script: update DC-Colos colo list
#
# The underlying data table (colo, name, lat, lon, cca2, region, city) had a
# row removed (CZX / Changzhou, China) which shifts all following rows up by
# one, and a brand new row was appended to the "Europe" block... actually to
# the very end of the existing data (just before Amman, Jordan) for
# BGI / Bridgetown, Barbados.
#
# Net effect, reproduced here with native Excel row operations so that all
# the in-between rows shift correctly without us having to rewrite every
# single cell by hand:
#   1. Delete the CZX (Changzhou, China) row entirely -> everything below
#      moves up by one row.
#   2. Insert a new blank row right before the "AMM" (Amman, Jordan) row
#      (which, after the deletion above, now sits one row higher than
#      before) and fill it in with the new Bridgetown, Barbados entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row containing the CZX colo code in column A and delete it.
$czxCell = $ws.Range("A1:A400").Find("CZX")
if ($czxCell -ne $null) {
    $czxRow = $czxCell.Row
    $ws.Rows($czxRow).Delete()
}

# Locate the row containing the AMM colo code in column A (Amman, Jordan).
# After the deletion above it has shifted up by one row already.
$ammCell = $ws.Range("A1:A400").Find("AMM")
$ammRow = $ammCell.Row

# Insert a new row right before the Amman row and populate it with the new
# Bridgetown, Barbados entry.
$ws.Rows($ammRow).Insert()

$newRow = $ammRow
$ws.Cells.Item($newRow, 1).Value = "BGI"
$ws.Cells.Item($newRow, 2).Value = "Bridgetown, Barbados"
$ws.Cells.Item($newRow, 3).Value = 13.103562
$ws.Cells.Item($newRow, 4).Value = -59.603226
$ws.Cells.Item($newRow, 5).Value = "BB"
$ws.Cells.Item($newRow, 6).Value = "North America"
$ws.Cells.Item($newRow, 7).Value = "Bridgetown"

# Match the formatting used by the other data rows: column A carries the
# bordered/bold/centered formatting (same as the rest of the table), while
# the other columns use the plain default formatting. Copy the format from
# the row above (which already has the correct look) rather than trying to
# reconstruct it cell-by-cell.
$ws.Cells.Item($newRow - 1, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
